$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated GDP per Capita values for existing years 1820-2010 (rows 2-192)
$gdpValues = @(
    "2031",
    "2104",
    "2115",
    "2085",
    "2117",
    "2107",
    "2110",
    "2150",
    "2163",
    "2110",
    "2120",
    "2101",
    "2158",
    "2130",
    "2227",
    "2195",
    "2179",
    "2219",
    "2216",
    "2224",
    "2276",
    "2254",
    "2243",
    "2353",
    "2434",
    "2472",
    "2501",
    "2466",
    "2570",
    "2695",
    "2817",
    "2649",
    "2713",
    "2697",
    "2679",
    "2938",
    "2735",
    "2745",
    "2691",
    "2834",
    "2775",
    "2785",
    "2836",
    "2979",
    "2922",
    "2989",
    "2963",
    "2933",
    "2954",
    "3092",
    "3193",
    "3177",
    "3327",
    "3279",
    "3341",
    "3366",
    "3395",
    "3261",
    "3351",
    "3425",
    "3476",
    "3480",
    "3571",
    "3665",
    "3642",
    "3625",
    "3724",
    "3818",
    "3808",
    "3826",
    "4022",
    "4073",
    "4141",
    "4191",
    "4235",
    "4415",
    "4521",
    "4564",
    "4575",
    "4705",
    "4809",
    "4948",
    "5007",
    "5244",
    "5302",
    "5333",
    "5423",
    "5557",
    "5662",
    "5807",
    "5906",
    "6148",
    "6076",
    "6236",
    "6551",
    "6022",
    "6202",
    "5765",
    "5514",
    "6153",
    "6363",
    "6099",
    "6641",
    "7267",
    "7218",
    "6978",
    "7329",
    "7425",
    "7627",
    "8089",
    "8513",
    "8542",
    "8239",
    "8434",
    "8611",
    "8735",
    "8886",
    "9035",
    "9185",
    "9553",
    "8155",
    "7291",
    "7379",
    "8097",
    "8835",
    "8075",
    "9208",
    "9620",
    "9776",
    "10351",
    "11067",
    "11056",
    "11086",
    "11623",
    "11749",
    "11788",
    "11858",
    "12696",
    "12903",
    "13767",
    "14046",
    "14843",
    "15537",
    "15513",
    "16832",
    "17459",
    "17789",
    "18230",
    "18868",
    "19974",
    "20221",
    "20617",
    "21579",
    "22228",
    "21919",
    "21712",
    "23059",
    "23360",
    "23632",
    "24409",
    "24272",
    "24063",
    "24807",
    "25450",
    "26581",
    "27710",
    "28681",
    "28728",
    "29049",
    "29108",
    "29412",
    "29915.4476167737",
    "30573.5733215054",
    "30649.5379873869",
    "32359.2962361468",
    "33356.1413655756",
    "34325.4425056305",
    "35499.0633285501",
    "36361.2294471395",
    "37523.4241694142",
    "39021.1775521356",
    "39425.8652231536",
    "39709.3692887464",
    "39983.1433535389",
    "41178.562625218",
    "42264.6302503897",
    "44025.4826886383",
    "44481.4687719655",
    "44246.3986882078",
    "42090.1732883123",
    "42932.3971449613"
)

$firstRow = 2
$lastRow = 192
$existingRange = $ws.Range("E$firstRow`:E$lastRow")
$existingRange.NumberFormat = "@"
for ($i = 0; $i -lt $gdpValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $gdpValues[$i]
}
$existingRange.ClearFormats()

# New rows for years 2011-2016 (rows 193-198)
$newValues = @(
    "43575",
    "43510",
    "43733",
    "44241",
    "44635",
    "44836"
)

$newFirstRow = 193
$newLastRow = 198
$newRange = $ws.Range("E$newFirstRow`:E$newLastRow")
$newRange.NumberFormat = "@"
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $newFirstRow + $i
    $year = 2011 + $i
    $ws.Range("A$row").Value = 208
    $ws.Range("B$row").Value = "Denmark"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $year
    $ws.Range("E$row").Value = $newValues[$i]
}
$newRange.ClearFormats()
